$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D) values that look numeric: stage them as text
# formulas in a scratch column, then paste-special (values-only) onto the
# target cells. A direct .Value assignment would let Excel re-interpret
# strings like "318.65" as numbers, which would lose formatting such as
# trailing zeros and flip the stored cell type from text to number.
$ws.Range("AA4").Formula = '="0.999"'
$ws.Range("AA5").Formula = '="318.65"'
$ws.Range("AA6").Formula = '="103.89"'
$ws.Range("AA7").Formula = '="0.516"'
$ws.Range("AA4:AA7").Copy()
$ws.Range("D4:D7").PasteSpecial(-4163)
$ws.Range("AA4:AA7").ClearContents()

$ws.Range("AA9").Formula = '="0.528"'
$ws.Range("AA10").Formula = '="35.84"'
$ws.Range("AA11").Formula = '="0.0804"'
$ws.Range("AA9:AA11").Copy()
$ws.Range("D9:D11").PasteSpecial(-4163)
$ws.Range("AA9:AA11").ClearContents()

$ws.Range("AA13").Formula = '="18.44"'
$ws.Range("AA14").Formula = '="7.00"'
$ws.Range("AA13:AA14").Copy()
$ws.Range("D13:D14").PasteSpecial(-4163)
$ws.Range("AA13:AA14").ClearContents()

$ws.Range("AA17").Formula = '="0.833"'
$ws.Range("AA17:AA17").Copy()
$ws.Range("D17:D17").PasteSpecial(-4163)
$ws.Range("AA17:AA17").ClearContents()

$ws.Range("AA19").Formula = '="12.33"'
$ws.Range("AA20").Formula = '="6.39"'
$ws.Range("AA19:AA20").Copy()
$ws.Range("D19:D20").PasteSpecial(-4163)
$ws.Range("AA19:AA20").ClearContents()

$ws.Range("AA22").Formula = '="69.09"'
$ws.Range("AA23").Formula = '="243.61"'
$ws.Range("AA22:AA23").Copy()
$ws.Range("D22:D23").PasteSpecial(-4163)
$ws.Range("AA22:AA23").ClearContents()

$ws.Range("AA25").Formula = '="2.51"'
$ws.Range("AA25:AA25").Copy()
$ws.Range("D25:D25").PasteSpecial(-4163)
$ws.Range("AA25:AA25").ClearContents()

$ws.Range("AA27").Formula = '="25.41"'
$ws.Range("AA28").Formula = '="2.24"'
$ws.Range("AA29").Formula = '="9.53"'
$ws.Range("AA30").Formula = '="33.56"'
$ws.Range("AA31").Formula = '="49.38"'
$ws.Range("AA32").Formula = '="0.127"'
$ws.Range("AA33").Formula = '="20.36"'
$ws.Range("AA34").Formula = '="5.23"'
$ws.Range("AA27:AA34").Copy()
$ws.Range("D27:D34").PasteSpecial(-4163)
$ws.Range("AA27:AA34").ClearContents()

$ws.Range("AA36").Formula = '="0.0767"'
$ws.Range("AA37").Formula = '="1.90"'
$ws.Range("AA38").Formula = '="4.50"'
$ws.Range("AA39").Formula = '="2.86"'
$ws.Range("AA40").Formula = '="124.74"'
$ws.Range("AA36:AA40").Copy()
$ws.Range("D36:D40").PasteSpecial(-4163)
$ws.Range("AA36:AA40").ClearContents()

$ws.Range("AA43").Formula = '="21.23"'
$ws.Range("AA44").Formula = '="0.0291"'
$ws.Range("AA43:AA44").Copy()
$ws.Range("D43:D44").PasteSpecial(-4163)
$ws.Range("AA43:AA44").ClearContents()

$ws.Range("AA46").Formula = '="2.94"'
$ws.Range("AA47").Formula = '="2.09"'
$ws.Range("AA48").Formula = '="9.25"'
$ws.Range("AA49").Formula = '="1.77"'
$ws.Range("AA50").Formula = '="76.03"'
$ws.Range("AA51").Formula = '="53.71"'
$ws.Range("AA46:AA51").Copy()
$ws.Range("D46:D51").PasteSpecial(-4163)
$ws.Range("AA46:AA51").ClearContents()

$excel.CutCopyMode = 0

# --- Update "Price" column (D) values that are not numeric-looking text
# (contain thousand separators as dots, or special glyphs) - plain .Value
# assignment keeps these as text automatically.
$ws.Range("D2").Value = "45.205.46"
$ws.Range("D3").Value = "2.438.63"
$ws.Range("D15").Value = "2.820.74"
$ws.Range("D16").Value = "2.440.66"
$ws.Range("D18").Value = "45.078.20"
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("D45").Value = "1.940.67"

# --- Update "Volume(1h)" column (E) values (always plain text: padded
# percentages), safe to assign directly.
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("E6").Value = "  +7.79%  "
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +9.39%  "
$ws.Range("E10").Value = "  +3.34%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").Value = "  +3.52%  "
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  +5.07%  "
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("E32").Value = "  +11.57%  "
$ws.Range("E33").Value = "  +13.51%  "
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +3.95%  "
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("E38").Value = "  +4.25%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  +6.55%  "
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +16.29%  "
$ws.Range("E50").Value = "  +6.19%  "
$ws.Range("E51").Value = "  +3.34%  "
